$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.720.77"
$ws.Range("E2").Value = "  -1.07%  "
$ws.Range("D3").Value = "2.984.55"
$ws.Range("E3").Value = "  -3.13%  "
$ws.Range("E4").Value = "  -0.41%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "497.91"
$ws.Range("E5").Value = "  -3.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.32"
$ws.Range("E6").Value = "  +4.49%  "
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.425"
$ws.Range("E8").Value = "  -2.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.21"
$ws.Range("E9").Value = "  +1.12%  "
$ws.Range("E10").Value = "  +1.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.349"
$ws.Range("E11").Value = "  -4.01%  "
$ws.Range("E12").Value = "  -0.49%  "
$ws.Range("D13").Value = "3.495.17"
$ws.Range("E13").Value = "  -3.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.12"
$ws.Range("E14").Value = "  +2.81%  "
$ws.Range("D15").Value = "56.691.72"
$ws.Range("E15").Value = "  -1.72%  "
$ws.Range("E16").Value = "  +1.98%  "
$ws.Range("D17").Value = "2.989.19"
$ws.Range("E17").Value = "  -3.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.67"
$ws.Range("E18").Value = "  +1.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.31"
$ws.Range("E19").Value = "  -2.74%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.76"
$ws.Range("E20").Value = "  +1.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "326.29"
$ws.Range("E21").Value = "  -2.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.468"
$ws.Range("E23").Value = "  -6.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.93"
$ws.Range("E24").Value = "  -5.89%  "
$ws.Range("E25").Value = "  -0.45%  "
$ws.Range("E26").Value = "  -2.01%  "
$ws.Range("D27").Value = "0.0₃0886"
$ws.Range("E27").Value = "  -0.32%  "
$ws.Range("E28").Value = "  -0.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.37"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.82"
$ws.Range("E30").Value = "  +2.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.73"
$ws.Range("E31").Value = "  -4.67%  "
$ws.Range("E32").Value = "  -5.95%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.25"
$ws.Range("E33").Value = "  -2.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "155.13"
$ws.Range("E34").Value = "  -1.48%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.41"
$ws.Range("E35").Value = "  -5.32%  "
$ws.Range("E36").Value = "  -3.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.55"
$ws.Range("E37").Value = "  -7.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0671"
$ws.Range("E38").Value = "  +1.12%  "
$ws.Range("B39").Value = "RenzoRestakedETH"
$ws.Range("C39").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D39").Value = "3.017.59"
$ws.Range("E39").Value = "  -3.55%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.80"
$ws.Range("E40").Value = "  +0.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.38"
$ws.Range("E41").Value = "  -8.91%  "
$ws.Range("E42").Value = "  -0.63%  "
$ws.Range("E43").Value = "  -5.56%  "
$ws.Range("D44").Value = "2.230.08"
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.992"
$ws.Range("E45").Value = "  -4.86%  "
$ws.Range("E46").Value = "  +0.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.54"
$ws.Range("E47").Value = "  -6.99%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.92"
$ws.Range("E48").Value = "  +11.68%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0235"
$ws.Range("E49").Value = "  +3.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.71"
$ws.Range("E50").Value = "  -4.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.90"
$ws.Range("E51").Value = "  -5.02%  "
